$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 940.9091
$ws.Range("I18").Value = 910
$ws.Range("K18").Value = 910
$ws.Range("M18").Value = -626
$ws.Range("H80").Value = 1322.2222
$ws.Range("J80").Value = 1100
$ws.Range("L80").Value = 3300
$ws.Range("N80").Value = -5296
$ws.Range("H83").Value = 1322.2222
$ws.Range("J83").Value = 1100
$ws.Range("L83").Value = 9900
$ws.Range("N83").Value = -19884
$ws.Range("H88").Value = 1701.5
$ws.Range("I88").Value = 1701.5
$ws.Range("K88").Value = 1701.5
$ws.Range("M88").Value = -1295.5
$ws.Range("H91").Value = 1701.5
$ws.Range("I91").Value = 1701.5
$ws.Range("K91").Value = 1701.5
$ws.Range("M91").Value = -297.5
$ws.Range("H105").Value = 20000
$ws.Range("J105").Value = 20000
$ws.Range("L105").Value = 20000
$ws.Range("N105").Value = -26988
$ws.Range("H111").Value = 659.125
$ws.Range("I111").Value = 395.5
$ws.Range("K111").Value = 1186.5
$ws.Range("M111").Value = 1880.5
$ws.Range("H138").Value = 2170.3333
$ws.Range("I138").Value = 707.4
$ws.Range("K138").Value = 2122.2
$ws.Range("M138").Value = 3017.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 250
$ws.Range("J4").Value = 250
$ws.Range("L4").Value = 250
$ws.Range("N4").Value = -482
$ws.Range("H6").Value = 599.5
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -127
$ws.Range("N12").ClearContents()
$ws.Range("H45").Value = 1216.1666
$ws.Range("I45").Value = 1132.8889
$ws.Range("K45").Value = 1132.8889
$ws.Range("M45").Value = -755.8888999999999
$ws.Range("H74").Value = 2200
$ws.Range("I74").Value = 2200
$ws.Range("K74").Value = 2200
$ws.Range("M74").Value = -1326
$ws.Range("H77").Value = 2200
$ws.Range("I77").Value = 2200
$ws.Range("K77").Value = 11000
$ws.Range("M77").Value = -6632
$ws.Range("H97").Value = 2392.7144
$ws.Range("I97").Value = 2039.8
$ws.Range("K97").Value = 2039.8
$ws.Range("M97").Value = -1543.8
$ws.Range("H110").Value = 1017.4
$ws.Range("I110").Value = 853.8570999999999
$ws.Range("K110").Value = 853.8570999999999
$ws.Range("M110").Value = 1191.1429
$ws.Range("H132").Value = 7882.2
$ws.Range("I132").Value = 7882.2
$ws.Range("K132").Value = 23646.6
$ws.Range("M132").Value = -21116.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 49999.5
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 49999.5
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 49999.5
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -50345.5
$ws.Range("H94").Value = 2312.3125
$ws.Range("I94").Value = 1833.1111
$ws.Range("K94").Value = 1833.1111
$ws.Range("M94").Value = -1382.1111

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1742.5454
$ws.Range("I58").Value = 1742.5454
$ws.Range("K58").Value = 1742.5454
$ws.Range("M58").Value = -1539.5454
$ws.Range("H99").Value = 2501000
$ws.Range("I99").Value = 2501000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2501000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2499502
$ws.Range("N99").ClearContents()
$ws.Range("H111").Value = 63250
$ws.Range("J111").Value = 63250
$ws.Range("L111").Value = 63250
$ws.Range("N111").Value = -71430
$ws.Range("H118").Value = 76000
$ws.Range("J118").Value = 76000
$ws.Range("L118").Value = 76000
$ws.Range("N118").Value = -79314
$ws.Range("H126").Value = 2501000
$ws.Range("I126").Value = 2501000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7503000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7500530
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 1742.5454
$ws.Range("I136").Value = 1742.5454
$ws.Range("K136").Value = 5227.6362
$ws.Range("M136").Value = -2677.6362

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 26415.5
$ws.Range("I32").Value = 7053.8887
$ws.Range("K32").Value = 21161.6661
$ws.Range("M32").Value = -20878.6661
$ws.Range("H46").Value = 1247.5
$ws.Range("I46").Value = 460
$ws.Range("J46").Value = 1810
$ws.Range("K46").Value = 1380
$ws.Range("L46").Value = 5430
$ws.Range("M46").Value = -1289
$ws.Range("N46").Value = -5612
$ws.Range("H97").Value = 191.23077
$ws.Range("I97").Value = 185.2
$ws.Range("J97").Value = 211.33333
$ws.Range("K97").Value = 555.5999999999999
$ws.Range("L97").Value = 633.99999
$ws.Range("M97").Value = -59.59999999999991
$ws.Range("N97").Value = -1625.99999
$ws.Range("H108").Value = 27
$ws.Range("I108").Value = 27
$ws.Range("K108").Value = 81
$ws.Range("M108").Value = 2799
$ws.Range("H111").Value = 150
$ws.Range("I111").Value = 150
$ws.Range("K111").Value = 450
$ws.Range("M111").Value = 2617
$ws.Range("H112").Value = 450
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H126").Value = 830
$ws.Range("I126").Value = 830
$ws.Range("K126").Value = 2490
$ws.Range("M126").Value = 2450
$ws.Range("H129").Value = 4749.75
$ws.Range("J129").Value = 4749.75
$ws.Range("L129").Value = 14249.25
$ws.Range("N129").Value = -24249.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 190
$ws.Range("I2").Value = 190
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 190
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -77
$ws.Range("N2").ClearContents()
$ws.Range("H102").Value = 24549.883
$ws.Range("J102").Value = 2014
$ws.Range("L102").Value = 2014
$ws.Range("N102").Value = -5258
$ws.Range("H122").Value = 2482.5
$ws.Range("I122").Value = 2354.5557
$ws.Range("J122").Value = 2866.3333
$ws.Range("K122").Value = 7063.6671
$ws.Range("L122").Value = 8598.999899999999
$ws.Range("M122").Value = -4613.6671
$ws.Range("N122").Value = -13498.9999
$ws.Range("H126").Value = 9999.5
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 976.3333
$ws.Range("I132").Value = 694.6
$ws.Range("K132").Value = 2083.8
$ws.Range("M132").Value = 446.1999999999998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7499
$ws.Range("I16").Value = 4999
$ws.Range("J16").Value = 9999
$ws.Range("K16").Value = 4999
$ws.Range("L16").Value = 9999
$ws.Range("M16").Value = -4829
$ws.Range("N16").Value = -10339
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H20").Value = 17000
$ws.Range("J20").Value = 17000
$ws.Range("L20").Value = 17000
$ws.Range("N20").Value = -17452
$ws.Range("H31").Value = 3464.3333
$ws.Range("J31").Value = 3844.2307
$ws.Range("L31").Value = 3844.2307
$ws.Range("N31").Value = -4340.2307
$ws.Range("H40").Value = 33000.777
$ws.Range("I40").Value = 23667.166
$ws.Range("K40").Value = 23667.166
$ws.Range("M40").Value = -23531.166
$ws.Range("H46").Value = 3733
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 4279.6
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 4279.6
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -4655.6
$ws.Range("H55").Value = 485.5
$ws.Range("I55").Value = 405.33334
$ws.Range("K55").Value = 405.33334
$ws.Range("M55").Value = -232.33334
$ws.Range("H61").Value = 2649.5
$ws.Range("I61").Value = 2649.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2649.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2447.5
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 2649.5
$ws.Range("I113").Value = 2649.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2649.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -479.5
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 2348.5
$ws.Range("I132").Value = 1964.6666
$ws.Range("K132").Value = 5893.9998
$ws.Range("M132").Value = -3363.9998
$ws.Range("H136").Value = 4931.1665
$ws.Range("I136").Value = 4757.4
$ws.Range("K136").Value = 14272.2
$ws.Range("M136").Value = -11722.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("K32").Value = 3000
$ws.Range("M32").Value = -2683
$ws.Range("H111").Value = 52666.668
$ws.Range("J111").Value = 52666.668
$ws.Range("L111").Value = 52666.668
$ws.Range("N111").Value = -60846.668
$ws.Range("H126").Value = 2008.091
$ws.Range("I126").Value = 1622.1111
$ws.Range("K126").Value = 4866.3333
$ws.Range("M126").Value = -2396.3333
$ws.Range("H136").Value = 7845.1665
$ws.Range("I136").Value = 7845.1665
$ws.Range("K136").Value = 23535.4995
$ws.Range("M136").Value = -20985.4995
